$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOM3066")

# Helper: set a cell's value as literal TEXT (not auto-converted to date/number)
# by round-tripping through a scratch cell formatted as Text, then Copy/PasteSpecial(values),
# which preserves the destination cell's existing style (so s="2"/s="3" survive untouched).
function Set-TextValue($cell, [string]$text) {
    $scratch = $ws.Cells.Item(1, 50)
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
    $scratch.Clear()
}

# 1) Ativação: date text 01/01/1996 -> 01/01/2022 (columns B and C)
Set-TextValue $ws.Cells.Item(8,2) '01/01/2022'
Set-TextValue $ws.Cells.Item(8,3) '01/01/2022'

# 2) Objetivos text (row 10)
$ws.Cells.Item(10,2).Value = 'Fornecer aos estudantes uma visão abrangente e interdisciplinar sobre materiais compósitos, além de mostrar as especificidades de cada matriz, sendo ela metálica, cerâmica ou polimérica. Ademais, deseja-se apresentar os fundamentos teóricos da mecânica de estruturas reforçadas e a partir de atividades práticas demostrar métodos de caracterização de materiais compósitos e como prepara-los.'
$ws.Cells.Item(10,3).Value = 'Fornecer aos estudantes uma visão abrangente e interdisciplinar sobre materiais compósitos, além de mostrar as especificidades de cada matriz, sendo ela metálica, cerâmica ou polimérica. Ademais, deseja-se apresentar os fundamentos teóricos da mecânica de estruturas reforçadas e a partir de atividades práticas demostrar métodos de caracterização de materiais compósitos e como prepara-los.'

# 3) Insert two new rows after row 13 for the two new responsible professors
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()
$ws.Cells.Item(14,2).Value = '1033242 - Fábio Herbst Florenzano'
$ws.Cells.Item(14,3).Value = '1033242 - Fábio Herbst Florenzano'
$ws.Cells.Item(15,2).Value = '1922320 - Sebastiao Ribeiro'
$ws.Cells.Item(15,3).Value = '1922320 - Sebastiao Ribeiro'

# 4) Programa resumido text (now row 16 after the insert)
$ws.Cells.Item(16,2).Value = '1.Introduçâo 2. Conceitos básicos sobre materiais compósitos, suas matrizes e seus processo de fabricação 3. Tipos de reforços 4. Compósitos nanoestruturados, naturais e híbridos 5. Mecânica da estrutura reforçada 6. Atividade prática'
$ws.Cells.Item(16,3).Value = '1.Introduçâo 2. Conceitos básicos sobre materiais compósitos, suas matrizes e seus processo de fabricação 3. Tipos de reforços 4. Compósitos nanoestruturados, naturais e híbridos 5. Mecânica da estrutura reforçada 6. Atividade prática'

# 5) Programa text (now row 18)
$ws.Cells.Item(18,2).Value = '1. Conceitos básicos sobre materiais compósitos: compósitos de matriz metálica (CMM), compósitos de matriz cerâmicos (CMC) e compósitos de matriz polimérica (CMP) e nanocompósitos. 2. Tipos de Reforços: Reforços particulados, fibras curtas, fibras longas, mantas, tecidos e preformas. 3. Conceitos de Interface4. Compósitos de matriz metálica: características e processos de fabricação. 5. Compósitos de matriz cerâmica: características e processos de fabricação. 6. Compósitos de matriz polimérica: matrizes termoplásticas e termorrígidas, características físicas e químicas e processos de fabricação. 7. Compósitos nanoestruturados. 8. Compósitos Naturais. 9. Compósitos Híbridos 10. Mecânica de estruturas reforçadas. Conteúdo prático: 1. Caracterização e análise de compósitos de matriz metálica. 2. Preparação e caracterização de compósitos de matriz polimérica.(Sugestão: Considerar substituir essa parte prática pela realização do PBL descrito no item 3) 3. Visita a empresa produtora de compósitos e aulas especiais e/ou palestras com professores/pesquisadores convidados'
$ws.Cells.Item(18,3).Value = '1. Conceitos básicos sobre materiais compósitos: compósitos de matriz metálica (CMM), compósitos de matriz cerâmicos (CMC) e compósitos de matriz polimérica (CMP) e nanocompósitos. 2. Tipos de Reforços: Reforços particulados, fibras curtas, fibras longas, mantas, tecidos e preformas. 3. Conceitos de Interface4. Compósitos de matriz metálica: características e processos de fabricação. 5. Compósitos de matriz cerâmica: características e processos de fabricação. 6. Compósitos de matriz polimérica: matrizes termoplásticas e termorrígidas, características físicas e químicas e processos de fabricação. 7. Compósitos nanoestruturados. 8. Compósitos Naturais. 9. Compósitos Híbridos 10. Mecânica de estruturas reforçadas. Conteúdo prático: 1. Caracterização e análise de compósitos de matriz metálica. 2. Preparação e caracterização de compósitos de matriz polimérica.(Sugestão: Considerar substituir essa parte prática pela realização do PBL descrito no item 3) 3. Visita a empresa produtora de compósitos e aulas especiais e/ou palestras com professores/pesquisadores convidados'

# 6) Método text (now row 21)
$ws.Cells.Item(21,2).Value = 'De acordo com a atual ementa da disciplina propõe-se o uso de uma nova metodologia de ensino com o intuito de abordar o conteúdo de forma mais prática e contextualizada para que o aluno consiga relacionar os conhecimentos teóricos vistos em sala de aula com as outras disciplinas do curso. Assim, avaliação do aluno será feita através de uma prova escrita e por uma apresentação final com base nas atividades práticas desenvolvidas.'
$ws.Cells.Item(21,3).Value = 'De acordo com a atual ementa da disciplina propõe-se o uso de uma nova metodologia de ensino com o intuito de abordar o conteúdo de forma mais prática e contextualizada para que o aluno consiga relacionar os conhecimentos teóricos vistos em sala de aula com as outras disciplinas do curso. Assim, avaliação do aluno será feita através de uma prova escrita e por uma apresentação final com base nas atividades práticas desenvolvidas.'

# 7) Critério text (now row 22)
$ws.Cells.Item(22,2).Value = 'A nota final será calculada como descrita a seguir: NF= (0,4*Avaliação escrita + 0,6 *Apresentação final)'
$ws.Cells.Item(22,3).Value = 'A nota final será calculada como descrita a seguir: NF= (0,4*Avaliação escrita + 0,6 *Apresentação final)'

# 8) Norma de recuperação text (now row 23)
$ws.Cells.Item(23,2).Value = 'Devido a cunho prático da disciplina não haverá recuperação.'
$ws.Cells.Item(23,3).Value = 'Devido a cunho prático da disciplina não haverá recuperação.'

# 9) Bibliografia text (now row 24)
$ws.Cells.Item(24,2).Value = '1. REZENDE, M. C.; COSTA, M. L.; BOTELHO, E. C. Compósitos estruturais: tecnologia e prática. São Paulo: Artliber, 2011. 396p. 2 MALLICK, P.K. Composites Engineering Handbook. New York: Marcel Dekker, 1997. 3. MATTHEWS, F.L. & RAWLINGS, R.D. Composite Materials: Engineering and Science. London: Chapman & Hall, 1994. 4. OBRAZTSOV, I.F. Mechanics of Composites. Moscow: MIR Publishers, 1982. 5. JONES R. Mechanics of Composite Materials. New York: McGraw-Hill, 1975. 6. UPADHYAYA, G.S. Sintered Metal-Ceramic Composites. Elsevier, 1984. 7. HARPER, C. A. Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill, 1992. 8. GOLDSTEIN, A.N. Handbook of Nanophase Materials. CRC Press, 1997. 9. DRESSELHAUS, M.S. Graphite Fibers and Filaments. New York: Springer-Verlag, 1988.'
$ws.Cells.Item(24,3).Value = '1. REZENDE, M. C.; COSTA, M. L.; BOTELHO, E. C. Compósitos estruturais: tecnologia e prática. São Paulo: Artliber, 2011. 396p. 2 MALLICK, P.K. Composites Engineering Handbook. New York: Marcel Dekker, 1997. 3. MATTHEWS, F.L. & RAWLINGS, R.D. Composite Materials: Engineering and Science. London: Chapman & Hall, 1994. 4. OBRAZTSOV, I.F. Mechanics of Composites. Moscow: MIR Publishers, 1982. 5. JONES R. Mechanics of Composite Materials. New York: McGraw-Hill, 1975. 6. UPADHYAYA, G.S. Sintered Metal-Ceramic Composites. Elsevier, 1984. 7. HARPER, C. A. Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill, 1992. 8. GOLDSTEIN, A.N. Handbook of Nanophase Materials. CRC Press, 1997. 9. DRESSELHAUS, M.S. Graphite Fibers and Filaments. New York: Springer-Verlag, 1988.'

